$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "38.778.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.093.49"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.46%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("E6").Value = "  +0.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.55"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.55%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +1.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0845"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.35%  "
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("B12").Value = "Chainlink"
$ws.Range("C12").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.53%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.406.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.814"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.48"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.085.77"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "38.731.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0842"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.48"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.78%  "
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.33%  "
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "171.37"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.140"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.78%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +11.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.53%  "
$ws.Range("E31").Value = "  +0.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.42"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.53"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.51%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.71"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0611"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.11%  "
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("E37").Value = "  -0.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.61"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.31%  "
$ws.Range("E40").Value = "  -0.99%  "
$ws.Range("E41").Value = "  +5.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.539.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.82%  "
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.79%  "
$ws.Range("B45").Value = "Cronos"
$ws.Range("C45").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0917"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.74"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.93%  "
$ws.Range("E47").Value = "  +1.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("E49").Value = "  +0.86%  "
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.294.71"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.72%  "
